# Updated cryptos list on Fri Sep  6 05:40:21 UTC 2024 with GitHub Actions
#
# Applies the per-row Price (D) / Volume(1h) (E) refresh, plus the two
# row swaps (USDe<->Aptos at rows 32/33, Aave<->RenderToken at rows 43/44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as TEXT, never letting Excel reinterpret
# a numeric-looking string (e.g. "506.46") as a float and mangle it with
# binary floating point noise (506.45999999999998). Plain non-numeric
# strings (multi-dot prices like "56.301.38", subscript prices, etc.) are
# assigned directly since Excel already stores them as text.
function Set-TextValue {
    param($range, [string]$text)
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

# --- Row 2 : Bitcoin ---
Set-TextValue $ws.Range("D2") '56.301.38'
$ws.Range("E2").Value = '  -1.35%  '

# --- Row 3 : Ethereum ---
Set-TextValue $ws.Range("D3") '2.370.92'
$ws.Range("E3").Value = '  -1.19%  '

# --- Row 4 : TetherUSD ---
$ws.Range("E4").Value = '  -0.02%  '

# --- Row 5 : BNB ---
Set-TextValue $ws.Range("D5") '506.46'

# --- Row 6 : Solana ---
Set-TextValue $ws.Range("D6") '129.89'
$ws.Range("E6").Value = '  -2.32%  '

# --- Row 7 : USDC ---
Set-TextValue $ws.Range("D7") '0.998'
$ws.Range("E7").Value = '  +0.02%  '

# --- Row 8 : XRP ---
$ws.Range("E8").Value = '  -1.72%  '

# --- Row 9 : LidoStakedEther ---
Set-TextValue $ws.Range("D9") '2.380.50'
$ws.Range("E9").Value = '  -1.30%  '

# --- Row 10 : Dogecoin ---
Set-TextValue $ws.Range("D10") '0.0985'
$ws.Range("E10").Value = '  +1.45%  '

# --- Row 11 : TRON ---
$ws.Range("E11").Value = '  -0.10%  '

# --- Row 12 : Toncoin ---
Set-TextValue $ws.Range("D12") '4.87'
$ws.Range("E12").Value = '  +6.43%  '

# --- Row 13 : Cardano ---
$ws.Range("E13").Value = '  +1.89%  '

# --- Row 14 : WrappedliquidstakedEther2.0 ---
Set-TextValue $ws.Range("D14") '2.792.50'
$ws.Range("E14").Value = '  -1.28%  '

# --- Row 15 : WrappedBTC ---
Set-TextValue $ws.Range("D15") '56.264.26'
$ws.Range("E15").Value = '  -1.28%  '

# --- Row 16 : Avalanche ---
Set-TextValue $ws.Range("D16") '21.50'
$ws.Range("E16").Value = '  -1.90%  '

# --- Row 17 : ShibaInu ---
Set-TextValue $ws.Range("D17") '0.0000133'
$ws.Range("E17").Value = '  -0.99%  '

# --- Row 18 : WrappedEther ---
Set-TextValue $ws.Range("D18") '2.377.30'

# --- Row 19 : Chainlink ---
Set-TextValue $ws.Range("D19") '10.02'
$ws.Range("E19").Value = '  -2.09%  '

# --- Row 20 : Polkadot ---
$ws.Range("E20").Value = '  -0.09%  '

# --- Row 21 : BitcoinCash ---
Set-TextValue $ws.Range("D21") '309.00'
$ws.Range("E21").Value = '  -0.34%  '

# --- Row 22 : Uniswap ---
Set-TextValue $ws.Range("D22") '6.29'
$ws.Range("E22").Value = '  -0.71%  '

# --- Row 23 : Dai ---
$ws.Range("E23").Value = '  -0.02%  '

# --- Row 24 : Litecoin ---
Set-TextValue $ws.Range("D24") '65.90'
$ws.Range("E24").Value = '  +1.45%  '

# --- Row 25 : Binance-PegBSC-USD ---
Set-TextValue $ws.Range("D25") '0.997'
$ws.Range("E25").Value = '  -0.01%  '

# --- Row 26 : Polygon ---
$ws.Range("E26").Value = '  -1.60%  '

# --- Row 27 : Kaspa ---
Set-TextValue $ws.Range("D27") '0.149'
$ws.Range("E27").Value = '  -2.92%  '

# --- Row 28 : InternetComputer(DFINITY) ---
Set-TextValue $ws.Range("D28") '7.21'
$ws.Range("E28").Value = '  -3.09%  '

# --- Row 29 : Monero ---
Set-TextValue $ws.Range("D29") '173.02'
$ws.Range("E29").Value = '  -0.53%  '

# --- Row 30 : PEPE ---
Set-TextValue $ws.Range("D30") '0.0₃0710'
$ws.Range("E30").Value = '  -1.93%  '

# --- Rows 32/33 : USDe and Aptos swap positions ---
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D32") '5.83'
$ws.Range("E32").Value = '  -2.06%  '

$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D33") '0.999'
$ws.Range("E33").Value = '  +0.02%  '

# --- Row 34 : FirstDigitalUSD ---
Set-TextValue $ws.Range("D34") '0.996'
$ws.Range("E34").Value = '  +0.08%  '

# --- Row 35 : Fetch.AI ---
$ws.Range("E35").Value = '  -5.02%  '

# --- Row 36 : EthereumClassic ---
Set-TextValue $ws.Range("D36") '17.60'
$ws.Range("E36").Value = '  -1.98%  '

# --- Row 38 : NEARProtocol ---
$ws.Range("E38").Value = '  -4.02%  '

# --- Row 39 : SuiNetwork ---
Set-TextValue $ws.Range("D39") '0.828'
$ws.Range("E39").Value = '  +1.91%  '

# --- Row 40 : OKB ---
Set-TextValue $ws.Range("D40") '36.36'
$ws.Range("E40").Value = '  -1.03%  '

# --- Row 41 : Stacks ---
Set-TextValue $ws.Range("D41") '1.39'
$ws.Range("E41").Value = '  -3.90%  '

# --- Row 42 : Filecoin ---
$ws.Range("E42").Value = '  +0.24%  '

# --- Rows 43/44 : Aave and RenderToken swap positions ---
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D43") '4.85'
$ws.Range("E43").Value = '  -3.22%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D44") '126.00'
$ws.Range("E44").Value = '  -5.45%  '

# --- Row 45 : Mantle ---
$ws.Range("E45").Value = '  -0.93%  '

# --- Row 46 : Stellar ---
$ws.Range("E46").Value = '  -1.39%  '

# --- Row 47 : Bittensor ---
Set-TextValue $ws.Range("D47") '238.05'
$ws.Range("E47").Value = '  -6.01%  '

# --- Row 48 : Hedera ---
$ws.Range("E48").Value = '  -1.81%  '

# --- Row 49 : VeChain ---
Set-TextValue $ws.Range("D49") '0.0207'
$ws.Range("E49").Value = '  -1.91%  '

# --- Row 50 : InjectiveProtocol ---
Set-TextValue $ws.Range("D50") '16.96'
$ws.Range("E50").Value = '  -1.48%  '

# --- Row 51 : BitgetToken ---
$ws.Range("E51").Value = '  -0.16%  '
